$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old B4 value (RegistrationData.InformationSource.Patient)
$ws.Range("B4").ClearContents()

# Insert a new row at 23, shifting current rows 23-24 down to 24-25
$ws.Rows("23").Insert()

# Set the new B23 to the value that used to live in B4
$ws.Range("B23").Value = "RegistrationData.InformationSource.Patient"

# Update the active selection to match the target state
$ws.Range("B23").Select()
